# Updated cryptos list (price/volume refresh) matching the GitHub Actions commit.
# Price cells that look like plain numbers (single '.' decimal, e.g. "331.66")
# are written with a leading apostrophe so Excel keeps them as literal text
# (preserving trailing zeros / multi-dot "thousands" formatting like
# "24.954.53") instead of silently coercing them into floating point values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.954.53"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "1.673.69"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").Value = "'331.66"
$ws.Range("E5").Value = "  +7.92%  "
$ws.Range("D6").Value = "'0.9993"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "'0.3658"
$ws.Range("E7").Value = "  +1.45%  "
$ws.Range("D8").Value = "'47.25"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").Value = "'0.3231"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("E10").Value = "  +2.28%  "
$ws.Range("D11").Value = "'0.07148"
$ws.Range("E11").Value = "  +3.05%  "
$ws.Range("D12").Value = "'0.9999"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "'6.091"
$ws.Range("E13").Value = "  +3.72%  "
$ws.Range("D14").Value = "'19.64"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("D15").Value = "1.668.09"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "'6.656"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "'0.06551"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "'0.9994"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("E20").Value = "  +3.49%  "
$ws.Range("E21").Value = "  +1.76%  "
$ws.Range("D22").Value = "'5.915"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").Value = "'12.81"
$ws.Range("E23").Value = "  +2.60%  "
$ws.Range("D24").Value = "24.950.06"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("D25").Value = "'2.439"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").Value = "'2.396"
$ws.Range("E26").Value = "  +4.65%  "
$ws.Range("D27").Value = "'148.70"
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("D29").Value = "1.856.68"
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("D30").Value = "'125.71"
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").Value = "'1.186"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").Value = "'4.089"
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("D33").Value = "'5.797"
$ws.Range("E33").Value = "  +3.50%  "
$ws.Range("D34").Value = "'0.08487"
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("D35").Value = "'1.662"
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("D36").Value = "'12.32"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").Value = "'5.165"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("D39").Value = "'1.234"
$ws.Range("E39").Value = "  +2.76%  "
$ws.Range("D40").Value = "'0.02230"
$ws.Range("E40").Value = "  +2.25%  "
$ws.Range("D41").Value = "'0.2092"
$ws.Range("E41").Value = "  +2.43%  "
$ws.Range("D42").Value = "'8.245"
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("D43").Value = "'0.9989"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "'0.5969"
$ws.Range("E44").Value = "  +1.68%  "
$ws.Range("D45").Value = "'13.71"
$ws.Range("E45").Value = "  +8.48%  "
$ws.Range("D46").Value = "'3.849"
$ws.Range("E46").Value = "  +3.14%  "
$ws.Range("D47").Value = "'0.5733"
$ws.Range("E47").Value = "  +3.18%  "
$ws.Range("D48").Value = "'124.22"
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("E49").Value = "  +1.86%  "
$ws.Range("D50").Value = "'0.07003"
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("D51").Value = "'1.195"
$ws.Range("E51").Value = "  +4.29%  "
